# Weekly CompStat report update (NYPD 7th Precinct)
# - Bumps the report volume/number and week-ending dates in the header
# - Refreshes every Crime Complaints statistic cell (rows 15-31) with newly
#   collected counts / period-over-period percentage changes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# ---- Header: volume/number + reporting week dates ----
$ws.Range("A8").Value = "Volume 32   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# ---- Crime Complaints table (rows 15-31) ----
# A handful of cells flip between a numeric count/pct-change and the
# sentinel text placeholders "0" / "***.*" used when a value is not
# applicable. Use Range.Copy(Destination) from a same-shaped template
# cell so both the display format and the literal text/number land
# correctly, then overwrite with the final value where needed.

$ws.Range("M15").Value = 25
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 11.111111111111
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = 17.241379310344
$ws.Range("L16").Value = -6.422018348623
$ws.Range("M16").Value = -12.068965517241
$ws.Range("N16").Value = -84.684684684684
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -58.333333333333
$ws.Range("I17").Value = 175
$ws.Range("J17").Value = 197
$ws.Range("K17").Value = -11.167512690355
$ws.Range("L17").Value = -4.371584699453
$ws.Range("M17").Value = 60.550458715596
$ws.Range("N17").Value = -7.894736842105
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 350
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -12.903225806451
$ws.Range("L18").Value = -44.137931034482
$ws.Range("M18").Value = 3.846153846153
$ws.Range("N18").Value = -75
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 357
$ws.Range("J19").Value = 427
$ws.Range("K19").Value = -16.393442622950
$ws.Range("L19").Value = -26.694045174538
$ws.Range("M19").Value = 58.666666666666
$ws.Range("N19").Value = -6.052631578947
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -41.666666666666
$ws.Range("L20").Value = -42.857142857142
$ws.Range("M20").Value = -36.363636363636
$ws.Range("N20").Value = -91.222570532915
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -38.888888888888
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -25.842696629213
$ws.Range("I21").Value = 755
$ws.Range("J21").Value = 861
$ws.Range("K21").Value = -12.311265969802
$ws.Range("L21").Value = -23.194303153611
$ws.Range("M21").Value = 29.948364888123
$ws.Range("N21").Value = -60.325801366263
$ws.Range("D14").Copy($ws.Range("C22"))
$ws.Range("F22").Value = 1
$ws.Range("D14").Copy($ws.Range("G22"))
$ws.Range("N22").Copy($ws.Range("H22"))
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -87.5
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = -41.666666666666
$ws.Range("I23").Value = 142
$ws.Range("J23").Value = 183
$ws.Range("K23").Value = -22.404371584699
$ws.Range("L23").Value = 5.185185185185
$ws.Range("M23").Value = 24.561403508771
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = -19.166666666666
$ws.Range("I24").Value = 1073
$ws.Range("J24").Value = 1149
$ws.Range("K24").Value = -6.614447345517
$ws.Range("L24").Value = 4.174757281553
$ws.Range("M24").Value = 71.132376395534
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -21.428571428571
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 625
$ws.Range("J25").Value = 765
$ws.Range("K25").Value = -18.300653594771
$ws.Range("L25").Value = 8.131487889273
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 7.692307692307
$ws.Range("F26").Value = 41
$ws.Range("H26").Value = -4.651162790697
$ws.Range("I26").Value = 344
$ws.Range("J26").Value = 362
$ws.Range("K26").Value = -4.972375690607
$ws.Range("L26").Value = -6.775067750677
$ws.Range("M26").Value = 27.407407407407
$ws.Range("D16").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("L14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = 14.705882352941
$ws.Range("M29").Value = 14.285714285714
$ws.Range("M30").Value = 25
$ws.Range("D14").Copy($ws.Range("G31"))
$ws.Range("N22").Copy($ws.Range("H31"))
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = -30.769230769230
$ws.Range("L31").Value = 200
